# Bugfixed evaluation and simulated rt_data for components
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-9 (B:F) with corrected values ---
$ws.Cells.Item(2,2).Value = 0.08631603385188587
$ws.Cells.Item(2,3).Value = 0.2336348552713851
$ws.Cells.Item(2,4).Value = 0.06595001910092338
$ws.Cells.Item(2,5).Value = 0.2568073579571337
$ws.Cells.Item(2,6).Value = 0.2509970670398742

$ws.Cells.Item(3,2).Value = 0.1074344412111381
$ws.Cells.Item(3,3).Value = 0.2042603668832133
$ws.Cells.Item(3,4).Value = 0.06469647195393771
$ws.Cells.Item(3,5).Value = 0.2543550116548477
$ws.Cells.Item(3,6).Value = 0.2399663286002585

$ws.Cells.Item(4,2).Value = 0.1080779521229239
$ws.Cells.Item(4,3).Value = 0.3459062368002563
$ws.Cells.Item(4,4).Value = 0.1618177455352025
$ws.Cells.Item(4,5).Value = 0.4022657648063063
$ws.Cells.Item(4,6).Value = 0.4047044737270317

$ws.Cells.Item(5,2).Value = 0.1348047727392753
$ws.Cells.Item(5,3).Value = 0.2404736632640445
$ws.Cells.Item(5,4).Value = 0.1336825235075403
$ws.Cells.Item(5,5).Value = 0.3656262073587455
$ws.Cells.Item(5,6).Value = 0.3564564719985849

$ws.Cells.Item(6,2).Value = 0.1094716894729319
$ws.Cells.Item(6,3).Value = 0.3502298540383924
$ws.Cells.Item(6,4).Value = 0.1684591240845271
$ws.Cells.Item(6,5).Value = 0.4104377225408589
$ws.Cells.Item(6,6).Value = 0.416966656392023

# Row 7: values update and rank (G7) changes 7 -> 9
$ws.Cells.Item(7,2).Value = 0.05568622347126046
$ws.Cells.Item(7,3).Value = 0.2806577130811624
$ws.Cells.Item(7,4).Value = 0.1125884625785131
$ws.Cells.Item(7,5).Value = 0.3355420429372646
$ws.Cells.Item(7,6).Value = 0.3509607463531708
$ws.Cells.Item(7,7).Value = 9

# Row 8: values update and rank (G8) changes 4 -> 6
$ws.Cells.Item(8,2).Value = -0.06433122593588818
$ws.Cells.Item(8,3).Value = 0.3695289678539498
$ws.Cells.Item(8,4).Value = 0.164134063030832
$ws.Cells.Item(8,5).Value = 0.4051346233424539
$ws.Cells.Item(8,6).Value = 0.4381719613125665
$ws.Cells.Item(8,7).Value = 6

# Row 9: values update (F9 newly populated) and rank (G9) changes 1 -> 3
$ws.Cells.Item(9,2).Value = -0.06392244592618833
$ws.Cells.Item(9,3).Value = 0.1995501394904485
$ws.Cells.Item(9,4).Value = 0.05897833981413506
$ws.Cells.Item(9,5).Value = 0.2428545651498754
$ws.Cells.Item(9,6).Value = 0.2869466694029099
$ws.Cells.Item(9,7).Value = 3

# --- Add new row 10 (Q8) ---
$ws.Cells.Item(10,1).Value = "Q8"
$ws.Cells.Item(10,2).Value = -0.06520887812495521
$ws.Cells.Item(10,3).Value = 0.06520887812495521
$ws.Cells.Item(10,4).Value = 0.004252197786315262
$ws.Cells.Item(10,5).Value = 0.06520887812495521
$ws.Cells.Item(10,7).Value = 1

# Copy formatting (border/bold/center) from A9 onto the new label cell A10
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
